# Card1: add a new service-log event row (2025-12-08, by HOSSAM).
#
# The sheet's "empty" data cells were originally exported from a pandas
# DataFrame where missing values serialize as the literal text "nan".
# Re-saving the sheet after appending the new record restores that
# literal "nan" text into every still-blank data cell (columns D:O on
# most rows, B:O on the otherwise-empty row 16) while leaving the cells
# that already carry real data untouched. The freshly appended row 17
# keeps its unused columns (B:K) genuinely blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# --- restore literal "nan" placeholders into currently-blank cells ---
$ws.Range("D2:K2").Value  = "nan"
$ws.Range("N2").Value     = "nan"

$ws.Range("G3:K3").Value  = "nan"
$ws.Range("M3:O3").Value  = "nan"

$ws.Range("D4:K4").Value  = "nan"

$ws.Range("D5:K5").Value  = "nan"

$ws.Range("D6:K6").Value  = "nan"
$ws.Range("M6").Value     = "nan"

$ws.Range("D7:O7").Value  = "nan"

$ws.Range("D8").Value     = "nan"
$ws.Range("H8").Value     = "nan"
$ws.Range("J8:K8").Value  = "nan"
$ws.Range("M8:O8").Value  = "nan"

$ws.Range("E9:G9").Value  = "nan"
$ws.Range("I9").Value     = "nan"
$ws.Range("K9").Value     = "nan"
$ws.Range("M9:O9").Value  = "nan"

$ws.Range("E10").Value    = "nan"
$ws.Range("G10:J10").Value = "nan"
$ws.Range("M10:O10").Value = "nan"

$ws.Range("E11:F11").Value = "nan"
$ws.Range("H11:K11").Value = "nan"
$ws.Range("M11:O11").Value = "nan"

$ws.Range("D12:O12").Value = "nan"
$ws.Range("D13:O13").Value = "nan"
$ws.Range("D14:O14").Value = "nan"
$ws.Range("D15:O15").Value = "nan"

$ws.Range("B16:K16").Value = "nan"

# --- append the new event as row 17 ---
# Column A ("card") is stored as text everywhere else in the column, so
# force text formatting before entering the value to avoid Excel
# reinterpreting "1" as a number.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "1"

$ws.Range("L17").Value = "30/8/2025"
$ws.Range("M17").Value = "قطع سير كويلر مسنن دبل 700"
$ws.Range("N17").Value = "تم تغير سير  دوبل700(محمد نعيم)"
$ws.Range("O17").Value = "فني"
